# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# tracker sheet to the latest scraped values, mirroring the automated
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Both columns were originally authored as plain text (inlineStr) cells
# with no cell style (`s` attribute), including Price values that look
# numeric (e.g. "8.58", "171.10", "1.00"). Writing such a string straight
# into Range.Value would make Excel auto-convert it into a real number
# (losing formatting like trailing zeros, e.g. "8.60" -> 8.6), so those
# cells are written using a leading quote (forcing literal text entry)
# and then their style is reset back to Normal so the cell keeps no
# explicit style, just like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => updated Price (D) / Volume(1h) (E) values.
$updates = @(
    @{ Row = 2;  D = '60.011.18';  E = '  +0.08%  ' },
    @{ Row = 3;  D = '2.410.09';   E = '  -0.24%  ' },
    @{ Row = 4;                    E = '  -0.07%  ' },
    @{ Row = 5;  D = '554.58';     E = '  +0.54%  ' },
    @{ Row = 6;  D = '136.31';     E = '  -0.66%  ' },
    @{ Row = 7;                    E = '  -0.04%  ' },
    @{ Row = 8;                    E = '  +0.83%  ' },
    @{ Row = 9;                    E = '  -1.19%  ' },
    @{ Row = 10; D = '5.62';       E = '  -2.34%  ' },
    @{ Row = 11;                   E = '  -0.62%  ' },
    @{ Row = 12; D = '0.352';      E = '  -1.36%  ' },
    @{ Row = 13; D = '24.74';      E = '  -0.11%  ' },
    @{ Row = 14; D = '2.841.78';   E = '  -0.24%  ' },
    @{ Row = 15; D = '59.875.02';  E = '  -0.06%  ' },
    @{ Row = 16;                   E = '  -0.12%  ' },
    @{ Row = 17; D = '2.409.78';   E = '  -0.67%  ' },
    @{ Row = 18; D = '11.19';      E = '  -0.92%  ' },
    @{ Row = 19;                   E = '  +3.40%  ' },
    @{ Row = 20; D = '326.47';     E = '  -1.37%  ' },
    @{ Row = 21; D = '6.75';       E = '  +0.76%  ' },
    @{ Row = 22;                   E = '  -0.05%  ' },
    @{ Row = 23; D = '64.73';      E = '  -1.73%  ' },
    @{ Row = 24;                   E = '  +3.52%  ' },
    @{ Row = 25; D = '8.60';       E = '  +0.15%  ' },
    @{ Row = 26; D = '0.999';      E = '  -0.75%  ' },
    @{ Row = 27;                   E = '  +3.92%  ' },
    @{ Row = 28;                   E = '  +1.32%  ' },
    @{ Row = 29; D = '0.0₃0769';   E = '  -1.42%  ' },
    @{ Row = 30; D = '171.16';     E = '  +0.29%  ' },
    @{ Row = 31; D = '6.12';       E = '  -1.23%  ' },
    @{ Row = 32;                   E = '  +8.20%  ' },
    @{ Row = 33;                   E = '  -2.98%  ' },
    @{ Row = 34; D = '18.41';      E = '  -1.09%  ' },
    @{ Row = 35;                   E = '  +0.04%  ' },
    @{ Row = 36;                   E = '  +2.51%  ' },
    @{ Row = 38; D = '4.21';       E = '  +0.73%  ' },
    @{ Row = 39; D = '322.55';     E = '  +2.63%  ' },
    @{ Row = 40;                   E = '  -0.76%  ' },
    @{ Row = 41; D = '146.90';     E = '  +6.10%  ' },
    @{ Row = 42;                   E = '  -1.85%  ' },
    @{ Row = 43; D = '0.0963';     E = '  +0.09%  ' },
    @{ Row = 44; D = '19.82';      E = '  +2.72%  ' },
    @{ Row = 45;                   E = '  -0.87%  ' },
    @{ Row = 46; D = '0.576';      E = '  -0.60%  ' },
    @{ Row = 47;                   E = '  -1.63%  ' },
    @{ Row = 48;                   E = '  +0.02%  ' },
    @{ Row = 49;                   E = '  -1.37%  ' },
    @{ Row = 50; D = '4.67';       E = '  -0.63%  ' },
    @{ Row = 51;                   E = '  -2.11%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('D')) {
        $dCell = $ws.Range("D$row")
        $newVal = $u.D

        $looksNumeric = ($newVal -as [double]) -ne $null

        if ($looksNumeric) {
            # Force literal text so Excel doesn't coerce it to a number
            # and strip the leading quote's formatting afterwards so the
            # cell is left without any explicit style, matching the rest
            # of the sheet.
            $dCell.Value = "'" + $newVal
            $dCell.Style = "Normal"
        }
        else {
            $dCell.Value = $newVal
        }
    }

    if ($u.ContainsKey('E')) {
        $ws.Range("E$row").Value = $u.E
    }
}
